# Update the model-prediction coefficients in the "depth/temp slope
# comparison" table with the refreshed estimates (better model
# predictions) from the commit "add better model predictions to coef
# analysis".
#
# Table layout (row, column):
#   col 1 contrast | col 2 estimate | col 3 std.error | col 4 df |
#   col 5 statistic | col 6 adj.p.value
#
# Cell values are addressed directly via Table.Cell(row, col) rather
# than Find/Replace because several values (e.g. "0.018", "-0.006",
# "22") repeat across multiple cells in this small table, and a
# Range-scoped Find here still matches across the whole story -- so a
# blind Find/Replace would clobber cells that must stay untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: "Deep - Mid"  -- estimate/df unchanged
$t.Cell(2, 3).Range.Text = "0.005"    # std.error   0.018 -> 0.005
$t.Cell(2, 5).Range.Text = "-1.045"   # statistic  -0.348 -> -1.045
$t.Cell(2, 6).Range.Text = "0.557"    # adj.p.value 0.935 -> 0.557

# Row 3: "Deep - Shallow" -- df unchanged
$t.Cell(3, 2).Range.Text = "-0.019"   # estimate   -0.012 -> -0.019
$t.Cell(3, 3).Range.Text = "0.006"    # std.error   0.018 -> 0.006
$t.Cell(3, 5).Range.Text = "-3.322"   # statistic  -0.680 -> -3.322
$t.Cell(3, 6).Range.Text = "0.008"    # adj.p.value 0.777 -> 0.008

# Row 4: "Mid - Shallow" -- df unchanged
$t.Cell(4, 2).Range.Text = "-0.013"   # estimate   -0.006 -> -0.013
$t.Cell(4, 3).Range.Text = "0.005"    # std.error   0.018 -> 0.005
$t.Cell(4, 5).Range.Text = "-2.460"   # statistic  -0.332 -> -2.460
$t.Cell(4, 6).Range.Text = "0.056"    # adj.p.value 0.941 -> 0.056
